# Horarios actualizados Linea 141 - 1147
# Refresh of the scraped bus-schedule snapshot: new "Ultima actualizacion" timestamp
# (03:08:51 -> 04:04:36) and a brand new set of scraped rows per sheet.

$wb = $excel.ActiveWorkbook

$newTs = "04:04:36"

# ---------------------------------------------------------------------------
# Sheet 1: LP1912  (all departures)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: $newTs"
$ws1.Range("A3").Value = "Total filas: 7"

$rows1 = @(
    @($newTs, "04:45", "215A_EL PATO",  41,  "LP1912"),
    @($newTs, "04:53", "11_ETCHEVERRY", 49,  "LP1912"),
    @($newTs, "05:16", "17_ROMERO",     72,  "LP1912"),
    @($newTs, "05:22", "23_HERNANDEZ",  78,  "LP1912"),
    @($newTs, "05:34", "215B_EL PATO",  90,  "LP1912"),
    @($newTs, "05:46", "15_ABASTO",     102, "LP1912"),
    @($newTs, "05:54", "10_OLMOS",      110, "LP1912")
)

$r = 6
foreach ($row in $rows1) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $r++
}

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215  (subset filtered on route 215)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: $newTs"
$ws2.Range("A3").Value = "Total filas: 2"

$rows2 = @(
    @($newTs, "04:45", "215A_EL PATO", 41, "LP1912"),
    @($newTs, "05:34", "215B_EL PATO", 90, "LP1912")
)

$r = 6
foreach ($row in $rows2) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $ws2.Cells.Item($r, 5).Value = $row[4]
    $r++
}

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173  (now has its first scraped row, plus column headers)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: $newTs"
$ws3.Range("A3").Value = "Total filas: 1"

$ws3.Cells.Item(5, 1).Value = "Hora_Scrap"
$ws3.Cells.Item(5, 2).Value = "Hora_Llegada"
$ws3.Cells.Item(5, 3).Value = "Linea"
$ws3.Cells.Item(5, 4).Value = "Minutos"
$ws3.Cells.Item(5, 5).Value = "Parada"

$ws3.Cells.Item(6, 1).Value = $newTs
$ws3.Cells.Item(6, 2).Value = "05:44"
$ws3.Cells.Item(6, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(6, 4).Value = 100
$ws3.Cells.Item(6, 5).Value = "L6173"
